$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "links in report" textbox (the shape carrying the GitHub
# repo hyperlink) robustly by its shape Id rather than a hard index.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Id -eq 8) {
        $sh = $cand
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(3)
}

# --- Resize / reposition the shape (EMU -> points, with a tiny epsilon to
# beat the interop layer's point->EMU truncation) ---
$emuPerPt = 12700
$eps = 0.00001
$sh.Left   = (5241409 / $emuPerPt) + $eps
$sh.Top    = (6465900 / $emuPerPt) + $eps
$sh.Width  = (5345084 / $emuPerPt) + $eps
$sh.Height = (392100  / $emuPerPt) + $eps

# --- Update the hyperlink run text, collapsing the 4 runs into 1 ---
$tr = $sh.TextFrame.TextRange
$orig_len = $tr.Length
$part = $tr.Characters(1, $orig_len)
$part.Text = "https://github.com/MCI/Combating MCI using Carrier Sensing"

# --- Remove the now-orphaned trailing empty paragraph ---
$tr2 = $sh.TextFrame.TextRange
$tail = $tr2.Characters($tr2.Length, 2)
$tail.Delete()
